$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the row above into the new row first,
# so the date cell in column G keeps the existing date/time number format.
$ws.Range("A9:H9").Copy($ws.Range("A10:H10"))

# Now fill in the actual values for the new row.
$ws.Cells.Item(10, 1).Value = 9702.51
$ws.Cells.Item(10, 2).Value = 9758.1299999999992
$ws.Cells.Item(10, 3).Value = 307.87
$ws.Cells.Item(10, 4).Value = 306.13
$ws.Cells.Item(10, 5).Value = $false
$ws.Cells.Item(10, 6).Value = -0.56999999999999995
$ws.Cells.Item(10, 7).Value = 42612.672905092593
$ws.Cells.Item(10, 8).Value = $false
